$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "usdf"
$ws.Range("B10").Value = "sdf@sdf.com"

$ws.Range("A11").Value = "sdsdf"
$ws.Range("B11").Value = "sdfsd@sdf.ccc"
